$p = $ppt.ActivePresentation

# --- 1. Slide 5: change the table's style (tblPr/tableStyleId) ---
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{9493FBA0-22B8-4EC7-A79F-3A78F61E84A9}")
    }
}

# --- 2. Swap the presentation's theme colours to the default "Office" palette ---
# (the Integral/"Red Violet" colour scheme is replaced with the stock Office colours)
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

# Order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
